# Moving the LPJ-GUESS variables whose fields duplicate other already-available
# variables from the "basic ignored" block into this "missing identified"
# file, plus one extra variable coming from the step-2 file (#141).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 10 blank rows starting at row 577 (pushes the existing
# "expfe / expsi / expcalc / ppdiat / ppmisc / SImon / areacellg / prra /
# mrroLi / cltisccp / ..." block down by 10 rows, matching the target
# layout where that block now starts at row 587).
$ws.Rows("577:586").Insert()

# Row 576 (pre-existing, still blank) + the 8 freshly inserted rows
# 577-584 now carry the new content; rows 585-586 stay blank.
$ws.Range("C576").Value2 = "mrsll"
$ws.Range("E576").Value2 = "Available in LPJ-GUESS, but the field will be the same as mrsol because: No frozen fraction of water"
$ws.Range("F576").Value2 = "David Wårlind"

$ws.Range("C577").Value2 = "netAtmosLandCO2Flux"
$ws.Range("E577").Value2 = "Available in LPJ-GUESS, but the field will be the same as nbp in Lmon"
$ws.Range("F577").Value2 = "David Wårlind"

$ws.Range("C578").Value2 = "nep"
$ws.Range("E578").Value2 = "Available in LPJ-GUESS, but the field will be the same as nbp in Lmon"
$ws.Range("F578").Value2 = "David Wårlind"

$ws.Range("C579").Value2 = "fFire"
$ws.Range("E579").Value2 = "Available in LPJ-GUESS, but the field will be the same as fFireNat"
$ws.Range("F579").Value2 = "David Wårlind"

$ws.Range("C580").Value2 = "fHarvest"
$ws.Range("E580").Value2 = "Available in LPJ-GUESS, but the field will be the same as fHarvestToAtmos"
$ws.Range("F580").Value2 = "David Wårlind"

$ws.Range("C581").Value2 = "cCwd"
$ws.Range("E581").Value2 = "Available in LPJ-GUESS, but the field will be the same as cLitterCwd"
$ws.Range("F581").Value2 = "David Wårlind"

$ws.Range("C582").Value2 = "rGrowth"
$ws.Range("E582").Value2 = "Available in LPJ-GUESS, but the field will be the same as raOther"
$ws.Range("F582").Value2 = "David Wårlind"

$ws.Range("C583").Value2 = "rMaint"
$ws.Range("E583").Value2 = "Available in LPJ-GUESS, but the field will be the same as r*"
$ws.Range("F583").Value2 = "David Wårlind"

# Row 584 also gets a Table entry (column A) -- the extra variable coming
# from the "step 2" file -- and is a touch taller than the default rows.
$ws.Range("A584").Value2 = "Eday"
$ws.Range("C584").Value2 = "prCrop"
$ws.Range("E584").Value2 = "Available in LPJ-GUESS, but the field will be the same precipitation over crops as for the rest of the gridcell. Available in LPJ-GUESS, but the field will be the same precipitation over crops as for the rest of the gridcell."
$ws.Range("F584").Value2 = "David Wårlind"
$ws.Range("E584").Font.Name = "Cambria"
$ws.Rows(584).RowHeight = 15

# Two pre-existing rows (now shifted to 604/605 -- "albisccp" / "pctisccp")
# get a new comment font: 12pt, dark grey, "normal arial".
$ws.Range("E604").Font.Name = "normal arial"
$ws.Range("E604").Font.Size = 12
$ws.Range("E604").Font.Color = 3355443
$ws.Range("E604").Copy()
$ws.Range("E605").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Reflect where the editor was last looking/scrolled to.
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$win.ScrollRow = 562
$win.ScrollColumn = 2
$ws.Range("B585").Select()
